$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Tên lớp" (class name) column is inserted as column B on just the
# top two header/value rows, pushing the existing B:H content of those two
# rows one column to the right (to C:I). The literal text below mirrors
# the existing cell contents of the template, written right-to-left so
# nothing is clobbered before its old value has been relocated.

# Row 1 (headers): A1 stays put, B1..H1 shift to C1..I1.
$ws.Cells.Item(1, 9).Value = "Mã Giảng viên"          # I1 <- old H1
$ws.Cells.Item(1, 8).Value = "Giảng viên"              # H1 <- old G1
$ws.Cells.Item(1, 7).Value = "Khóa học"                # G1 <- old F1
$ws.Cells.Item(1, 6).Value = "Khoa"                    # F1 <- old E1
$ws.Cells.Item(1, 5).Value = "Loại hình đào tạo"       # E1 <- old D1
$ws.Cells.Item(1, 4).Value = "Bậc đào tạo"             # D1 <- old C1
$ws.Cells.Item(1, 3).Value = "Chuyên ngành"            # C1 <- old B1
$ws.Cells.Item(1, 2).Value = "Tên lớp"                 # B1 <- brand new header

# New I1 cell needs the same look as the rest of the header row (thin
# border + light fill); give it the closest achievable direct formatting.
$ws.Cells.Item(1, 9).Borders.LineStyle = 1
$ws.Cells.Item(1, 9).Borders.Color = 0
$ws.Cells.Item(1, 9).Interior.Color = 10086143

# Row 2 (values): A2 stays put, B2..H2 shift to C2..I2.
$ws.Cells.Item(2, 9).Value = 1024571                              # I2 <- old H2
$ws.Cells.Item(2, 8).Value = "Nguyễn Năm"                          # H2 <- old G2
$ws.Cells.Item(2, 7).Value = "2016-2020"                           # G2 <- old F2
$ws.Cells.Item(2, 6).Value = "CNTT"                                # F2 <- old E2
$ws.Cells.Item(2, 5).Value = "Tiên tiến"                           # E2 <- old D2
$ws.Cells.Item(2, 4).Value = "Đại học"                             # D2 <- old C2
$ws.Cells.Item(2, 3).Value = "Kỹ Thuật Phần Mềm"                   # C2 <- old B2
$ws.Cells.Item(2, 2).Value = "Đại học Kỹ thuật Phần mềm 12B CLC"   # B2 <- brand new value

# New I2 cell: same thin border treatment as the rest of row 2 (no fill).
$ws.Cells.Item(2, 9).Borders.LineStyle = 1
$ws.Cells.Item(2, 9).Borders.Color = 0

# Restore the current selection to C11, as recorded in the saved view.
$ws.Range("C11").Select()
